# Applies the "add three new rotation rows for 24/03/23" edit to the
# WaiterRotationApp Cycle sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source cells already carrying the three distinct cell styles used in the
# sheet (style 1 = bold/bordered date header, style 2 = green fill,
# style 3 = red fill, and plain/default formatting from C2).
$styleSources = @{
    0 = "C2"
    1 = "A2"
    2 = "B2"
    3 = "D2"
}

function Set-CellStyle($address, $styleId) {
    $src = $styleSources[$styleId]
    $ws.Range($src).Copy()
    $ws.Range($address).PasteSpecial(-4122)
}

# (cell address, style id, value)
$newCells = @(
    @("A3", 1, "24/03/23"), @("B3", 3, "Table1"), @("C3", 3, "Table2-3"), @("D3", 2, "13-14-15"),
    @("E3", 0, "SmokingSection"), @("F3", 0, "SmokingSection"), @("G3", 0, "Middle"), @("H3", 3, "NewSection"),
    @("I3", 0, "Middle"), @("J3", 2, "16-17-18"), @("K3", 3, "Table4-5"), @("L3", 3, "NewSection"),

    @("A4", 1, "24/03/23"), @("B4", 0, "SmokingSection"), @("C4", 2, "16-17-18"), @("D4", 3, "NewSection"),
    @("E4", 3, "Table1"), @("F4", 0, "Middle"), @("G4", 3, "NewSection"), @("H4", 3, "Table2-3"),
    @("I4", 3, "NewSection"), @("J4", 0, "Middle"), @("K4", 2, "13-14-15"), @("L4", 3, "Table4-5"),

    @("A5", 1, "24/03/23"), @("B5", 3, "Table2-3"), @("C5", 3, "Table4-5"), @("D5", 0, "Middle"),
    @("E5", 0, "Middle"), @("F5", 3, "NewSection"), @("G5", 2, "13-14-15"), @("H5", 0, "SmokingSection"),
    @("I5", 0, "SmokingSection"), @("J5", 3, "Table1"), @("K5", 3, "NewSection"), @("L5", 2, "16-17-18")
)

foreach ($cell in $newCells) {
    Set-CellStyle $cell[0] $cell[1]
}
foreach ($cell in $newCells) {
    $ws.Range($cell[0]).Value = $cell[2]
}

# Move the "next empty row" selection down to row 6, same way it pointed at
# row 3 right after row 2 held the last entered rotation.
$ws.Rows.Item(6).Select() | Out-Null
